$d = $word.ActiveDocument

# 1. Append " and allow spaces in names" as its own run right after the
#    run containing "Names of nodes with auto-assignment" (same paragraph).
$r1 = $d.Content
$r1.Find.Execute("Names of nodes with auto-assignment", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0) | Out-Null
$r1.Collapse(0)
$r1.InsertAfter(" and allow spaces in names")
# Touch formatting on the inserted text so it stays a distinct run instead
# of silently re-merging with the preceding, identically formatted run.
$r1.Bold = 1
$r1.Bold = 0

# 2. Insert a new bullet paragraph right after the paragraph that contains
#    "Cover by tests file handlers and solver", at the top list level.
$r2 = $d.Content
$r2.Find.Execute("Cover by tests file handlers and solver", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$coverPara = $r2.Paragraphs(1)
$newPara = $coverPara.Next()
$newPara.Range.Text = "Not needed every time ask to save net. In only cases when the net has been changed."
$newPara.Range.ListFormat.ListLevelNumber = 1
